# ECU Validar Solicitud de Búsqueda.docx — apply the two semantic edits:
#   1. Set page orientation to Portrait (explicit) on every section.
#   2. Collapse the multi-run "La validación..." paragraph into a single run.

$d = $word.ActiveDocument

# --- 1. Force explicit Portrait orientation on every section's page setup ---
foreach ($sec in $d.Sections) {
    $sec.PageSetup.Orientation = 0   # wdOrientPortrait
}

# --- 2. Merge the fragmented runs of the introductory paragraph into one run ---
$oldText = "La validación de una solicitud de búsqueda se realiza cuando la unidad orgánica registra la solicitud, para que después esta pase a manos del técnico, quien se encargará de validar si la solicitud contiene todos los datos necesarios para proceder con el envío al secretario general y, de darse el escenario más favorable, se apruebe y proceda la posterior búsqueda del documento solicitado."

$range = $d.Content
$found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $oldText, 2)

Write-Output "Paragraph merge found/replaced: $found"
